$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.211.40'
$ws.Range('E2').Value = '  -2.79%  '

$ws.Range('D3').Value = '3.520.16'
$ws.Range('E3').Value = '  -4.74%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.64%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.35'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.71%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.609'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.25%  '

$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.513.72'
$ws.Range('E8').Value = '  -4.73%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.189'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.10%  '

$ws.Range('E11').Value = '  -1.97%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.585'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.23%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.26'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.78%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000274'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.65%  '

$ws.Range('D15').Value = '4.088.12'
$ws.Range('E15').Value = '  -4.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.54'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -5.22%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '629.52'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -7.28%  '

$ws.Range('D18').Value = '3.533.85'
$ws.Range('E18').Value = '  -4.48%  '

$ws.Range('D19').Value = '69.194.56'
$ws.Range('E19').Value = '  -3.11%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.123'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.46'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.91%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.19'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.45%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.887'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.98%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.96'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -8.52%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.66'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.39%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.80'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.67%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.64'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.79%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.35'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -9.18%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.75'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.75%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.16'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.56'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.71%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.33'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.45%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.03'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -7.04%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '634.46'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +8.96%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.76'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.91%  '

$ws.Range('E37').Value = '  -5.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.46'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -15.51%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '57.36'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.33%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0454'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.17%  '

$ws.Range('E42').Value = '  -5.76%  '

$ws.Range('D43').Value = '3.381.05'
$ws.Range('E43').Value = '  -8.05%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.329'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.78%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '32.91'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -8.08%  '

$ws.Range('D46').Value = '0.0₃0696'
$ws.Range('E46').Value = '  -9.61%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.57'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.36%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.77'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -5.17%  '

$ws.Range('E49').Value = '  -2.41%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.73'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +14.68%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '132.00'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.41%  '

